# Add sail related and base data
# Appends four new dialog rows (44-47) to the defaultDialog sheet, each
# introducing a brand-new shared string in column A, "4" in column B and
# zeros in columns C-F, following the exact pattern of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: dialog_sail_without_enough_sailors
$ws.Range("A44").Value = "dialog_sail_without_enough_sailors"
$ws.Range("B44").Value = 4
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0

# Row 45: dialog_sail_without_fillup
$ws.Range("A45").Value = "dialog_sail_without_fillup"
$ws.Range("B45").Value = 4
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0

# Row 46: dialog_days_enough_to_sail
$ws.Range("A46").Value = "dialog_days_enough_to_sail"
$ws.Range("B46").Value = 4
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0

# Row 47: dialog_days_not_enough_to_sail
$ws.Range("A47").Value = "dialog_days_not_enough_to_sail"
$ws.Range("B47").Value = 4
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0

# Move the view/selection down to the newly appended last row, matching the
# scrolled viewport the author ended up with after adding the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("A47").Select()
